$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.865.23"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "3.540.77"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'611.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.27%  "
$ws.Range("D6").Value = "'184.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +4.38%  "
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").Value = "'53.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("D13").Value = "'9.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "4.101.68"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").Value = "'608.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.74%  "
$ws.Range("D16").Value = "69.904.23"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.570.44"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'18.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").Value = "'4.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("D24").Value = "'99.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.82%  "
$ws.Range("D25").Value = "'4.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("E27").Value = "  -5.57%  "
$ws.Range("E28").Value = "  +4.51%  "
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").Value = "'7.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.53%  "
$ws.Range("D31").Value = "'12.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("E34").Value = "  +18.56%  "
$ws.Range("D35").Value = "'3.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.97%  "
$ws.Range("D36").Value = "'533.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.49%  "
$ws.Range("E37").Value = "  -4.96%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "'37.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0779"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'3.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.39%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.539.80"
$ws.Range("E42").Value = "  +5.10%  "
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("D44").Value = "'0.0456"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("D45").Value = "'2.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.69%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.142"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.73%  "
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").Value = "'1.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("D51").Value = "'135.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.49%  "
